$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-11 Sunday" "2025-05-12 Monday"

Replace-Text "70×16=1120" "26×33=858"
Replace-Text "62×44=2728" "11×77=847"
Replace-Text "60×44=2640" "26×25=650"
Replace-Text "32×91=2912" "76×36=2736"
Replace-Text "92×64=5888" "55×68=3740"

Replace-Text "55×35=1925" "65×32=2080"
Replace-Text "53×88=4664" "11×21=231"
Replace-Text "66×64=4224" "55×30=1650"
Replace-Text "54×60=3240" "23×73=1679"
Replace-Text "38×72=2736" "11×49=539"

Replace-Text "42×79=3318" "77×75=5775"
Replace-Text "24×27=648" "71×41=2911"
Replace-Text "12×50=600" "25×77=1925"
Replace-Text "95×72=6840" "39×27=1053"
Replace-Text "45×75=3375" "32×52=1664"

Replace-Text "54×57=3078" "56×92=5152"
Replace-Text "45×91=4095" "13×32=416"
Replace-Text "28×78=2184" "32×80=2560"
Replace-Text "69×32=2208" "31×51=1581"
Replace-Text "36×30=1080" "47×69=3243"

Replace-Text "58×38=2204" "53×98=5194"
Replace-Text "90×55=4950" "26×86=2236"
Replace-Text "57×93=5301" "91×86=7826"
Replace-Text "12×34=408" "82×73=5986"
Replace-Text "63×86=5418" "82×82=6724"
